$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.741029
$ws.Range("H2").Value = 23.223087
$ws.Range("I2").Value = 0.4930486933812723
$ws.Range("J2").Value = 0.4930486933812723
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.451416666666667
$ws.Range("N2").Value = 7.35425
$ws.Range("O2").Value = 0.2191928499183569
$ws.Range("P2").Value = 0.2191928499183569
$ws.Range("Q2").Value = 18.97648750775
$ws.Range("R2").Value = 170.78838756975
$ws.Range("S2").Value = 0.1080727482507632
$ws.Range("T2").Value = 0.1080727482507632

$ws.Range("G3").Value = 7.741029
$ws.Range("H3").Value = 23.223087
$ws.Range("I3").Value = 0.4930486933812723
$ws.Range("J3").Value = 0.4930486933812723
$ws.Range("O3").Value = 0.4446889938320204
$ws.Range("P3").Value = 0.4446889938320204
$ws.Range("Q3").Value = 38.498678854855
$ws.Range("R3").Value = 346.488109693695
$ws.Range("S3").Value = 0.2192533273699103
$ws.Range("T3").Value = 0.2192533273699103

$ws.Range("G4").Value = 7.741029
$ws.Range("H4").Value = 23.223087
$ws.Range("I4").Value = 0.4930486933812723
$ws.Range("J4").Value = 0.4930486933812723
$ws.Range("O4").Value = 0.3361181562496228
$ws.Range("P4").Value = 0.3361181562496228
$ws.Range("Q4").Value = 29.09922470361
$ws.Range("R4").Value = 261.89302233249
$ws.Range("S4").Value = 0.1657226177605988
$ws.Range("T4").Value = 0.1657226177605988

$ws.Range("I5").Value = 0.0194007766416684
$ws.Range("J5").Value = 0.0194007766416684
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.451416666666667
$ws.Range("N5").Value = 7.35425
$ws.Range("O5").Value = 0.2191928499183569
$ws.Range("P5").Value = 0.2191928499183569
$ws.Range("Q5").Value = 0.7466982481111112
$ws.Range("R5").Value = 6.720284233000001
$ws.Range("S5").Value = 0.004252511522716785
$ws.Range("T5").Value = 0.004252511522716785

$ws.Range("I6").Value = 0.0194007766416684
$ws.Range("J6").Value = 0.0194007766416684
$ws.Range("O6").Value = 0.4446889938320204
$ws.Range("P6").Value = 0.4446889938320204
$ws.Range("S6").Value = 0.008627311844343285
$ws.Range("T6").Value = 0.008627311844343285

$ws.Range("I7").Value = 0.0194007766416684
$ws.Range("J7").Value = 0.0194007766416684
$ws.Range("O7").Value = 0.3361181562496228
$ws.Range("P7").Value = 0.3361181562496228
$ws.Range("S7").Value = 0.006520953274608331
$ws.Range("T7").Value = 0.006520953274608331

$ws.Range("G8").Value = 7.654706000000001
$ws.Range("I8").Value = 0.4875505299770593
$ws.Range("J8").Value = 0.4875505299770593
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.451416666666667
$ws.Range("N8").Value = 7.35425
$ws.Range("O8").Value = 0.2191928499183569
$ws.Range("P8").Value = 0.2191928499183569
$ws.Range("Q8").Value = 18.76487386683334
$ws.Range("R8").Value = 168.8838648015
$ws.Range("S8").Value = 0.1068675901448769
$ws.Range("T8").Value = 0.1068675901448769

$ws.Range("G9").Value = 7.654706000000001
$ws.Range("I9").Value = 0.4875505299770593
$ws.Range("J9").Value = 0.4875505299770593
$ws.Range("O9").Value = 0.4446889938320204
$ws.Range("P9").Value = 0.4446889938320204
$ws.Range("Q9").Value = 38.06936623313668
$ws.Range("S9").Value = 0.2168083546177668
$ws.Range("T9").Value = 0.2168083546177668

$ws.Range("G10").Value = 7.654706000000001
$ws.Range("I10").Value = 0.4875505299770593
$ws.Range("J10").Value = 0.4875505299770593
$ws.Range("O10").Value = 0.3361181562496228
$ws.Range("P10").Value = 0.3361181562496228
$ws.Range("Q10").Value = 28.77472877754001
$ws.Range("S10").Value = 0.1638745852144156
$ws.Range("T10").Value = 0.1638745852144156

